# Generate Report for Archive
# - Update "Status" text from "Ready for handoff" to "In Translation"
#   (shows up on Overview!E2:F3 and on zh-cn!C2:C3 / de-de!C2:C3)
# - Narrow the "Status" column width on Overview (E:F) and on the
#   zh-cn / de-de sheets (C) from 17.2159881591797 to 13.4101845877511

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Replace the status text wherever it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Narrow the status columns (closest value the pixel-quantized ColumnWidth
# model can land on for the target stored width of 13.4101845877511).
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
